# Commit: "update file with jgit"
# The only meaningful content change is cell E8 on the "Rules" sheet:
# its text changes from "Good Morning" to "GIT UPDATE". (The shared
# strings table reshuffling seen in the raw XML diff is a natural
# consequence of Excel dropping the now-unused "Good Morning" string
# and appending the new "GIT UPDATE" string — the engine handles that
# automatically.) The diff also shows E8 becoming the active/selected
# cell in the sheet view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E8").Value = "GIT UPDATE"
$ws.Range("E8").Select()
